$d = $word.ActiveDocument

# Replace "<id>p005r_a2</id>" with "<id>p005r_2</id>"
$d.Content.Find.Execute("<id>p005r_a2</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p005r_2</id>", 2)

# Replace "<id>p005v_a1</id>" with "<id>p005v_1</id>"
$d.Content.Find.Execute("<id>p005v_a1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p005v_1</id>", 2)
